# cambio de fracciones e historico
# Update the "Reporte de Formatos" sheet: refresh the reporting period (Q2 -> Q3 2022)
# and rename the responsible-area text; tweak a few display/formatting details to
# match how the workbook looked after the real edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Data updates in row 8 -------------------------------------------------
# Fecha de inicio del periodo que se informa
$ws.Range("B8").Value = 44743
# Fecha de término del periodo que se informa
$ws.Range("C8").Value = 44834
# Fecha de validación / Fecha de Actualización
$ws.Range("F8").Value = 44844
$ws.Range("G8").Value = 44844

# Área(s) responsable(s) ... text change (E8 keeps the plain-text style, only the
# text itself changes)
$ws.Range("E8").Value = "Departamento de Recursos Humanos (UPP)"

# --- 2. Row / column sizing ----------------------------------------------------
$ws.Rows.Item(3).RowHeight = 63.75
$ws.Rows.Item(8).RowHeight = 75

$ws.Columns.Item(4).ColumnWidth = 79.5
$ws.Columns.Item(6).ColumnWidth = 21.333333333333332
$ws.Columns.Item(7).ColumnWidth = 28
$ws.Columns.Item(8).ColumnWidth = 37

# --- 3. Border tweak on the hyperlink cell (D8): drop the left/top edges so
# only the right & bottom edges keep their thin border -------------------------
$d8 = $ws.Range("D8")
$d8.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft -> xlLineStyleNone
$d8.Borders.Item(8).LineStyle = -4142   # xlEdgeTop  -> xlLineStyleNone

# --- 4. Sheet view: where the window was scrolled/selected when saved ---------
$ws.Activate()
$ws.Range("G15").Select()

# --- 5. Print orientation ------------------------------------------------------
$ws.PageSetup.Orientation = 1   # xlPortrait

Write-Host "edit complete"
